$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 9
    3  = 9
    4  = 5
    5  = 3
    6  = 9
    7  = 6
    8  = 4
    9  = 6
    10 = 2
    11 = 6
    12 = 1
    13 = 3
    14 = 2
    15 = 4
    16 = 5
    17 = 8
    18 = 1
    19 = 7
    20 = 10
    21 = 6
    22 = 2
    23 = 10
    24 = 8
    25 = 3
    26 = 9
    27 = 4
    28 = 7
    29 = 4
    30 = 4
    31 = 11
    32 = 6
    33 = 6
    34 = 10
    35 = 2
    36 = 3
    37 = 4
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
